$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# NSF date should be present: set end_date for the current (ongoing) appointment row to "Present"
$ws.Range("C5").Value = "Present"

# Select C6 to mirror the post-edit active cell
$ws.Range("C6").Select()
